$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: correct-answer weight changes from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row: total marks recomputed (10 right * 5 = 50)
$ws.Range("B12").Value = 50

# Correct / total marks summary text
$ws.Range("E12").Value = "50/140"
